# Add a new review row (row 14) to the review database sheet, mirroring
# the formatting of the previous row (row 13), and wire up the two
# mailto hyperlinks for the email/recovery columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Values for the new row ---
$ws.Range("A14").Value = "com.singleton.strechy"
$ws.Range("B14").Value = "taxi game"
$ws.Range("C14").Value = "vikicrestina@gmail.com"
$ws.Range("D14").Value = "cristianjohn1222@gmail.com"
$ws.Range("E14").Value = "27/5/2019 15:59"
$ws.Range("F14").Value = "best taxi game ever – guaranteed!!"

# --- Hyperlinks for the email / recovery-email columns ---
$ws.Hyperlinks.Add($ws.Range("C14"), "mailto:vikicrestina@gmail.com", "", "", "vikicrestina@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D14"), "mailto:cristianjohn1222@gmail.com", "", "", "cristianjohn1222@gmail.com")

# --- Copy the formatting from row 13 onto row 14 (values/hyperlinks above are preserved) ---
$ws.Range("A13:F13").Copy()
$ws.Range("A14:F14").PasteSpecial(-4122)

# --- Match the selection left by the edit (active cell C14) ---
$ws.Range("C14").Select()
